# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric; force text format so Excel
# does not silently convert them to numbers (matches original inlineStr cells).
$textCells = @("D5", "D6", "D9", "D10", "D12", "D15", "D17", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.884.66"
$ws.Range("E2").Value = "  -3.47%  "
$ws.Range("D3").Value = "3.209.92"
$ws.Range("E3").Value = "  -4.72%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "534.93"
$ws.Range("E5").Value = "  -5.89%  "
$ws.Range("D6").Value = "135.06"
$ws.Range("E6").Value = "  -9.33%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.209.52"
$ws.Range("E8").Value = "  -4.75%  "
$ws.Range("D9").Value = "0.458"
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").Value = "7.53"
$ws.Range("E10").Value = "  -5.53%  "
$ws.Range("E11").Value = "  -6.73%  "
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  -5.77%  "
$ws.Range("D13").Value = "3.763.28"
$ws.Range("E13").Value = "  -4.70%  "
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "25.83"
$ws.Range("E15").Value = "  -7.84%  "
$ws.Range("D16").Value = "3.212.26"
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "0.0000157"
$ws.Range("E17").Value = "  -7.29%  "
$ws.Range("D18").Value = "58.886.45"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("D19").Value = "5.91"
$ws.Range("E19").Value = "  -6.93%  "
$ws.Range("D20").Value = "13.30"
$ws.Range("E20").Value = "  -8.11%  "
$ws.Range("D21").Value = "8.21"
$ws.Range("E21").Value = "  -8.06%  "
$ws.Range("D22").Value = "361.11"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "70.13"
$ws.Range("D25").Value = "0.519"
$ws.Range("E25").Value = "  -7.61%  "
$ws.Range("D26").Value = "3.339.77"
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("D27").Value = "0.171"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("E28").Value = "  -11.00%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "7.10"
$ws.Range("E30").Value = "  -5.03%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  -7.73%  "
$ws.Range("D33").Value = "7.06"
$ws.Range("E33").Value = "  -8.68%  "
$ws.Range("D34").Value = "21.74"
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  -7.73%  "
$ws.Range("D36").Value = "161.80"
$ws.Range("E36").Value = "  -5.02%  "
$ws.Range("D37").Value = "4.89"
$ws.Range("E37").Value = "  -9.08%  "
$ws.Range("D38").Value = "6.36"
$ws.Range("E38").Value = "  -6.86%  "
$ws.Range("E39").Value = "  -8.58%  "
$ws.Range("D40").Value = "26.22"
$ws.Range("E40").Value = "  -9.78%  "
$ws.Range("D41").Value = "0.0706"
$ws.Range("E41").Value = "  -6.69%  "
$ws.Range("D42").Value = "3.241.76"
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("D43").Value = "40.86"
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("D44").Value = "0.714"
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "4.02"
$ws.Range("E46").Value = "  -6.58%  "
$ws.Range("E47").Value = "  -7.02%  "
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "2.302.63"
$ws.Range("E49").Value = "  -7.62%  "
$ws.Range("D50").Value = "6.27"
$ws.Range("E50").Value = "  -6.36%  "
$ws.Range("D51").Value = "20.74"
$ws.Range("E51").Value = "  -8.33%  "

# Reset style on the forced-text cells back to Normal so no stray
# style index is left attached to the cell (matches original, unstyled cells).
foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}
